$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 137 (shifts existing rows 137:198 down to 138:199)
$ws.Rows("137:137").Insert()

# Populate the newly inserted row 137 with the new record
$ws.Range("A137").Value = 5
$ws.Range("B137").Value = "Macroferia Regional de Talca"
$ws.Range("C137").Value = "Maule"
$ws.Range("D137").Value = 44466
$ws.Range("E137").Value = 7
$ws.Range("F137").Value = 100112032
$ws.Range("G137").Value = "Zapallo italiano"
$ws.Range("H137").Value = "Sin especificar"
$ws.Range("I137").Value = "Primera"
$ws.Range("J137").Value = 300
$ws.Range("K137").Value = 10000
$ws.Range("L137").Value = 10000
$ws.Range("M137").Value = 10000
$ws.Range("N137").Value = "`$/caja 50 unidades"
$ws.Range("O137").Value = "Región de Arica y Parinacota"
$ws.Range("P137").Value = 200
$ws.Range("Q137").Value = 50
$ws.Range("R137").Value = "Hortaliza"
